# Adds the "known tech stack" parenthetical (bold+italic) right after each
# job-title run, and bumps the "Full Stack Developer" end date from
# "Present" to "April 2023".
#
# For each of the three experience entries we:
#   1. Find the exact title text and replace it with the same text plus a
#      trailing space (keeps it in the same run/formatting - italic only).
#   2. Insert the bold+italic tech-stack text right after that, then mark
#      just that inserted span Bold.
# Finally, the Full Stack Developer date range text is updated in place.

$d = $word.ActiveDocument

function Add-TechStack($TitleText, $TechStack) {
    $rng = $d.Content
    $found = $rng.Find.Execute($TitleText, $true, $false, $false, $false, $false, $true, 1, $false, ($TitleText + " "), 2)
    if (-not $found) {
        throw "Could not find title text: $TitleText"
    }

    $insPoint = $rng.End
    $insRng = $d.Range($insPoint, $insPoint)
    $insRng.InsertAfter($TechStack)

    $boldRng = $d.Range($insPoint, $insPoint + $TechStack.Length)
    $boldRng.Bold = 1
}

Add-TechStack "Full Stack Developer" "(Node.js, Python, GoLang, AWS, Flask)"
Add-TechStack "SDE intern" "(Augmented Reality, React.js, Unity3d, C#, AWS)"
Add-TechStack "Virtual Reality Development Intern" "( Virtual Reality, OculusVR, Unity3d, C#)"

# Update the Full Stack Developer role end date: "Present" -> "April 2023".
$dateRng = $d.Content
$dateFound = $dateRng.Find.Execute("May 2022 - Present", $true, $false, $false, $false, $false, $true, 1, $false, "May 2022 - April 2023", 2)
if (-not $dateFound) {
    throw "Could not find the Full Stack Developer date range"
}
